# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1488
$ws1.Range("F3").Value = 3138
$ws1.Range("F5").Value = 861
$ws1.Range("F6").Value = 296

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1488
$ws4.Range("F3").Value = 3138
$ws4.Range("F5").Value = 861
$ws4.Range("F7").Value = 296
